# Append newly scraped Lancers listings (2025-11-18 12:37:43 JST run),
# merged+sorted by priority score (column G, descending) with the previous
# run's rows, and bump every row's captured-at timestamp to the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clean slate: drop existing hyperlinks/relationships so they can be
# re-added below in final row order (F2..F11 -> rId1..rId10).
$ws.Hyperlinks.Delete()

# Clear old data rows (2:4); new data is written fresh below.
$ws.Range("A2:H4").ClearContents()

$rows = @(
    @{A='2025-11-18 12:37:43'; B='【急募】業種判定AIツールのGAS開発依頼'; C='システム開発'; D='5,000 円 ~ 10,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436501'; G=405; H='🔥AI,Ai ◆ツール,開発'},
    @{A='2025-11-18 12:37:43'; B='【謝礼あり】AIに興味のあるエンジニアの方へ|45分インタビュー(2,000円)協力お願いします'; C='システム開発'; D='1,000 ~ 5,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436391'; G=295; H='🔥AI,Ai'},
    @{A='2025-11-18 12:37:43'; B='【業務委託】Shopee価格調整ツールの開発(Googleスプレッドシート+GAS)'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436149'; G=128; H='◆ツール,開発'},
    @{A='2025-11-18 12:37:43'; B='初回 【急募】エンタメ型ガチャアプリのMVP開発|Web+iOS対応フルスタックエンジニア募集'; C='システム開発'; D='1,000,000 円 ~ 3,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436594'; G=100; H='◆開発 ◇アプリ'},
    @{A='2025-11-18 12:37:43'; B='初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5425629'; G=45; H='◇サイト'},
    @{A='2025-11-18 12:37:43'; B='【技術パートナー募集】リード獲得・育成システム構築'; C='システム開発'; D='100,000 円 ~ 200,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436021'; G=33; H=''},
    @{A='2025-11-18 12:37:43'; B='〖リモート可〗Delphiエンジニア募集'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5341051'; G=25; H=''},
    @{A='2025-11-18 12:37:43'; B='Blender担当講師募集(Roblox向け3Dモデリング/完全リモート/長期歓迎)'; C='システム開発'; D='100,000 円 ~ 200,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436476'; G=18; H=''},
    @{A='2025-11-18 12:37:43'; B='【Robloxクリエイター育成】講師募集!完全リモート可'; C='システム開発'; D='10,000 円 ~ 20,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436426'; G=10; H=''},
    @{A='2025-11-18 12:37:43'; B='【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え'; C='システム開発'; D='10,000 円 ~ 20,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5436248'; G=10; H=''}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    if ($row.H -ne "") {
        $ws.Cells.Item($r, 8).Value = $row.H
    }
    $r = $r + 1
}

# Column widths (B/D/H widened to fit the longer new titles/prices/skill tags).
$ws.Columns.Item(2).ColumnWidth = 54.16666666666667
$ws.Columns.Item(4).ColumnWidth = 31.16666666666667
$ws.Columns.Item(8).ColumnWidth = 15.16666666666667

# Stash the pre-existing "Hyperlink" cell style (from F2) in a scratch
# cell first -- Hyperlinks.Add() below re-stamps its own style variant
# on whatever cell it targets, so F2 itself cannot be used as the source
# again afterwards.
$ws.Range("F2").Copy() | Out-Null
$ws.Range("Z1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-create the URL hyperlinks for every data row (order F2..F11 ->
# rId1..rId10).
for ($row = 2; $row -le 11; $row++) {
    $url = $ws.Cells.Item($row, 6).Value()
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $url) | Out-Null
}

# Re-apply the clean "Hyperlink" style to every link cell, then drop the
# scratch cell.
$ws.Range("Z1").Copy() | Out-Null
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 6).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear() | Out-Null

